$d = $word.ActiveDocument

# Replace the author name and capture the resulting range (covers the
# newly-inserted "Hollis Veal" text).
$rng = $d.Content
$rng.Find.Execute("David Thomas", $true, $false, $false, $false, $false, $true, 1, $false, "Hollis Veal", 2)

# Word tracks the location of the most recent edit with the hidden
# "_GoBack" bookmark. Re-adding a bookmark with that name moves it here
# (and removes the old one automatically, since bookmark names are
# unique), matching what Word does after an edit like this.
$d.Bookmarks.Add("_GoBack", $rng)
